$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update F2 on the data sheet (refreshed panel_query_time)
$ws.Range("F2").Value = "2021-10-05 14:22:14.008843"

# Add new "metadata" sheet right after "data"
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Pityriasis rubra pilaris"
$meta.Range("C2").Value = 311
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.1"
$meta.Range("E2").Value = "2017-11-05T02:37:20.407928Z"
$meta.Range("F2").Value = "2021-10-05 14:22:14.005850"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/311/?format=json"

# Match the header styling (bold, bordered, centered) used on the "data" sheet
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# D2 must stay plain (unstyled) like B2/C2/E2/F2/G2 -- reset the number
# format applied above back to the sheet's default look
$meta.Range("B2").Copy()
$meta.Range("D2").PasteSpecial(-4122)  # xlPasteFormats
